# preparation publication 0.2.0
# - bump Version metadata value to 0.2.0
# - bump Date metadata value
# - insert a new "Jurisdiction" row (with value "iso:code:3166:FR") right
#   after the "Contact" row, pushing Description/Purpose/Copyright/Immutable
#   down by one row each

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple value updates -------------------------------------------------
$ws.Range("B3").Value = "0.2.0"
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# --- make room for the new "Jurisdiction" row -----------------------------
# Shift rows 11..14 down to 12..15 (bottom-up so we never overwrite data we
# still need to read). Formats are copied first (so a brand-new row 15
# picks up the same style as the row it is cloned from), then the actual
# cell values are written explicitly so empty source cells really end up
# empty in the destination.
for ($r = 14; $r -ge 11; $r--) {
    $dest = $r + 1
    $srcRange = "A" + $r + ":B" + $r
    $dstRange = "A" + $dest + ":B" + $dest

    $aVal = $ws.Range("A" + $r).Value2
    $bVal = $ws.Range("B" + $r).Value2

    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122)

    $ws.Range("A" + $dest).Value = $aVal
    if ($bVal -eq $null) {
        $ws.Range("B" + $dest).ClearContents()
    } else {
        $ws.Range("B" + $dest).Value = $bVal
    }
}

# --- write the new row ------------------------------------------------
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
